$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update phone numbers in existing rows (01, 02, 03) ---
# Rows 2 and 3: phone becomes plain text "0823456789" (leading apostrophe forces
# text storage so the leading zero survives; Excel strips the prefix marker itself).
$ws.Range("F2").Value = "'0823456789"
$ws.Range("F3").Value = "'0823456789"

# Row 4: phone becomes text that literally starts with an apostrophe character.
# A doubled leading apostrophe means "force text" (first quote) followed by a
# literal apostrophe (second quote) that is kept as real cell content.
$ws.Range("F4").Value = "''0823456789"

# --- Append new row 7 (test case 06) ---
$ws.Range("A7").Value = "'06"
$ws.Range("B7").Value = "Ngoc Vu"
$ws.Range("C7").Value = "207 Giai Phong"
$ws.Range("D7").Value = "Ha Noi"
$ws.Range("E7").Value = "'100000"
$ws.Range("F7").Value = "'0823456789"
$ws.Range("G7").Value = "Hmm. We couldn’t find your device’s location. Please enter the address manually."

# Row 7 formatting: wrap text on every cell like the other data rows, plus an
# explicit black font on the name/address/city cells (matches new style used
# for this row), and a taller row to fit the wrapped message.
$ws.Range("A7:G7").WrapText = $true
$ws.Range("B7:D7").Font.Color = 0
$ws.Rows.Item(7).RowHeight = 43.2

# --- Update the active selection shown in the sheet view ---
$ws.Range("E11").Select()
